$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (data rows only, skip the header) to text so numeric-looking
# strings (e.g. "310.34") are not silently converted to floating point
# numbers by Excel's type inference.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.898.29'
$ws.Range('E2').Value = '  -2.12%  '
$ws.Range('D3').Value = '1.831.21'
$ws.Range('E3').Value = '  -1.95%  '
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').Value = '310.34'
$ws.Range('E5').Value = '  -1.83%  '
$ws.Range('D6').Value = '1.005'
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('D7').Value = '0.4616'
$ws.Range('E7').Value = '  -1.03%  '
$ws.Range('D8').Value = '0.3651'
$ws.Range('E8').Value = '  -2.34%  '
$ws.Range('D9').Value = '0.07169'
$ws.Range('E9').Value = '  -3.04%  '
$ws.Range('D10').Value = '0.8780'
$ws.Range('E10').Value = '  -1.26%  '
$ws.Range('D11').Value = '0.07825'
$ws.Range('E11').Value = '  -1.77%  '
$ws.Range('D12').Value = '19.58'
$ws.Range('E12').Value = '  -2.44%  '
$ws.Range('D13').Value = '1.872.32'
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('D14').Value = '5.322'
$ws.Range('E14').Value = '  -2.08%  '
$ws.Range('D15').Value = '6.372'
$ws.Range('E15').Value = '  -3.64%  '
$ws.Range('D16').Value = '88.48'
$ws.Range('E16').Value = '  -4.65%  '
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').Value = '0.000008738'
$ws.Range('E18').Value = '  -2.52%  '
$ws.Range('D19').Value = '1.005'
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('D20').Value = '26.922.24'
$ws.Range('E20').Value = '  -2.15%  '
$ws.Range('D21').Value = '14.49'
$ws.Range('E21').Value = '  -3.07%  '
$ws.Range('D22').Value = '5.001'
$ws.Range('E22').Value = '  -3.59%  '
$ws.Range('E23').Value = '  -1.67%  '
$ws.Range('D24').Value = '1.972'
$ws.Range('E24').Value = '  +4.44%  '
$ws.Range('D25').Value = '150.68'
$ws.Range('D26').Value = '18.20'
$ws.Range('E26').Value = '  -2.06%  '
$ws.Range('D27').Value = '1.995'
$ws.Range('E27').Value = '  -4.82%  '
$ws.Range('D28').Value = '113.53'
$ws.Range('E28').Value = '  -3.35%  '
$ws.Range('D29').Value = '4.942'
$ws.Range('E29').Value = '  -4.47%  '
$ws.Range('D30').Value = '0.08824'
$ws.Range('E30').Value = '  -1.09%  '
$ws.Range('D31').Value = '3.103'
$ws.Range('E31').Value = '  +2.59%  '
$ws.Range('D32').Value = '0.7592'
$ws.Range('E32').Value = '  +0.71%  '
$ws.Range('E33').Value = '  -0.80%  '
$ws.Range('D34').Value = '1.136'
$ws.Range('E34').Value = '  -2.21%  '
$ws.Range('D35').Value = '2.652'
$ws.Range('E35').Value = '  -0.91%  '
$ws.Range('D36').Value = '1.090'
$ws.Range('E36').Value = '  +0.43%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '2.917'
$ws.Range('E38').Value = '  -2.52%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.05123'
$ws.Range('E39').Value = '  -3.53%  '
$ws.Range('D40').Value = '6.936'
$ws.Range('E40').Value = '  -3.79%  '
$ws.Range('D41').Value = '0.4973'
$ws.Range('E41').Value = '  -5.14%  '
$ws.Range('D42').Value = '0.1594'
$ws.Range('E42').Value = '  -3.25%  '
$ws.Range('D43').Value = '8.359'
$ws.Range('E43').Value = '  -0.10%  '
$ws.Range('D44').Value = '10.25'
$ws.Range('E44').Value = '  -0.89%  '
$ws.Range('D45').Value = '0.4664'
$ws.Range('E45').Value = '  -4.98%  '
$ws.Range('E46').Value = '  +0.27%  '
$ws.Range('D47').Value = '102.49'
$ws.Range('E47').Value = '  -1.41%  '
$ws.Range('D48').Value = '1.608'
$ws.Range('E48').Value = '  -3.65%  '
$ws.Range('D49').Value = '0.06096'
$ws.Range('E49').Value = '  -2.65%  '
$ws.Range('D50').Value = '64.53'
$ws.Range('E50').Value = '  -2.31%  '
$ws.Range('E51').Value = '  -2.66%  '

# Restore column D formatting/style so no stray number format leaks into the
# saved styles (cells keep their original default style).
$ws.Range("D2:D51").NumberFormat = "General"
$ws.Range("D2:D51").Style = "Normal"
